$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.8506941505955009
$ws0.Range("C2").Value = -0.5449213174164297
$ws0.Range("B3").Value = 0.6663400042645358
$ws0.Range("C3").Value = -0.9902754239070344
$ws0.Range("B4").Value = 1.284140086668387
$ws0.Range("C4").Value = -0.8720229218888095

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.449479897556458
$ws1.Range("C2").Value = 0.1271650696554648
$ws1.Range("B3").Value = 1.223815471612321
$ws1.Range("C3").Value = -0.01243671900592642
$ws1.Range("B4").Value = -0.8800887217874873
$ws1.Range("C4").Value = -0.6697415577546783
